$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 75
$ws.Range("I4").Value = 77
$ws.Range("K4").Value = 77
$ws.Range("M4").Value = 37

$ws.Range("H86").Value = 5914.643
$ws.Range("I86").Value = 5551.2856
$ws.Range("K86").Value = 5551.2856
$ws.Range("M86").Value = -4428.2856

$ws.Range("H89").Value = 5914.643
$ws.Range("I89").Value = 5551.2856
$ws.Range("K89").Value = 27756.428
$ws.Range("M89").Value = -22140.428

$ws.Range("H92").Value = 799.8
$ws.Range("I92").Value = 612.25
$ws.Range("K92").Value = 612.25
$ws.Range("M92").Value = 635.75

$ws.Range("H129").Value = 5078.875
$ws.Range("J129").Value = 4203.9165
$ws.Range("L129").Value = 12611.7495
$ws.Range("N129").Value = -22611.7495

$ws.Range("H137").Value = 1424396.6
$ws.Range("I137").Value = 17865.223
$ws.Range("J137").Value = 4800072
$ws.Range("K137").Value = 53595.66900000001
$ws.Range("L137").Value = 14400216
$ws.Range("M137").Value = -51045.66900000001
$ws.Range("N137").Value = -14405316

$ws.Range("H138").Value = 3684.7563
$ws.Range("J138").Value = 4418.0728
$ws.Range("L138").Value = 13254.2184
$ws.Range("N138").Value = -23534.2184

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4737.5386
$ws.Range("I61").Value = 4459.9
$ws.Range("J61").Value = 5663
$ws.Range("K61").Value = 4459.9
$ws.Range("L61").Value = 5663
$ws.Range("M61").Value = -4247.9
$ws.Range("N61").Value = -6087

$ws.Range("H74").Value = 22248050
$ws.Range("I74").Value = 27729.395
$ws.Range("J74").Value = 142872640
$ws.Range("K74").Value = 27729.395
$ws.Range("L74").Value = 142872640
$ws.Range("M74").Value = -26855.395
$ws.Range("N74").Value = -142874388

$ws.Range("H77").Value = 22248050
$ws.Range("I77").Value = 27729.395
$ws.Range("J77").Value = 142872640
$ws.Range("K77").Value = 138646.975
$ws.Range("L77").Value = 714363200
$ws.Range("M77").Value = -134278.975
$ws.Range("N77").Value = -714371936

$ws.Range("H132").Value = 5805.6665
$ws.Range("I132").Value = 5976.1113
$ws.Range("J132").Value = 5550
$ws.Range("K132").Value = 17928.3339
$ws.Range("L132").Value = 16650
$ws.Range("M132").Value = -15398.3339
$ws.Range("N132").Value = -21710

$ws.Range("H136").Value = 4737.5386
$ws.Range("I136").Value = 4459.9
$ws.Range("J136").Value = 5663
$ws.Range("K136").Value = 13379.7
$ws.Range("L136").Value = 16989
$ws.Range("M136").Value = -10829.7
$ws.Range("N136").Value = -22089

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 7714.409
$ws.Range("I20").Value = 6837.294
$ws.Range("K20").Value = 6837.294
$ws.Range("M20").Value = -6590.294

$ws.Range("H99").Value = 2571.8262
$ws.Range("I99").Value = 912
$ws.Range("K99").Value = 912
$ws.Range("M99").Value = 586

$ws.Range("H107").Value = 1825.7084
$ws.Range("I107").Value = 1312.1666
$ws.Range("J107").Value = 3366.3333
$ws.Range("K107").Value = 1312.1666
$ws.Range("L107").Value = 3366.3333
$ws.Range("M107").Value = 607.8334
$ws.Range("N107").Value = -7206.3333

$ws.Range("H134").Value = 4153.6665
$ws.Range("I134").Value = 3755.5
$ws.Range("K134").Value = 11266.5
$ws.Range("M134").Value = -8731.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12208.7
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 12208.7
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 12208.7
$ws.Range("N31").Value = -12798.7
$ws.Range("M31").ClearContents()

$ws.Range("H34").Value = 12208.7
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 12208.7
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 12208.7
$ws.Range("N34").Value = -12612.7
$ws.Range("M34").ClearContents()

$ws.Range("H58").Value = 2674.0588
$ws.Range("I58").Value = 2496.3333
$ws.Range("K58").Value = 2496.3333
$ws.Range("M58").Value = -2293.3333

$ws.Range("H99").Value = 3241.6572
$ws.Range("I99").Value = 3502.2856
$ws.Range("J99").Value = 2199.1428
$ws.Range("K99").Value = 3502.2856
$ws.Range("L99").Value = 2199.1428
$ws.Range("M99").Value = -2004.2856
$ws.Range("N99").Value = -5195.1428

$ws.Range("H107").Value = 831.17645
$ws.Range("I107").Value = 949.6429000000001
$ws.Range("J107").Value = 278.33334
$ws.Range("K107").Value = 949.6429000000001
$ws.Range("L107").Value = 278.33334
$ws.Range("M107").Value = 970.3570999999999
$ws.Range("N107").Value = -4118.33334

$ws.Range("H122").Value = 1245.2858
$ws.Range("J122").Value = 1367.5
$ws.Range("L122").Value = 4102.5
$ws.Range("N122").Value = -9002.5

$ws.Range("H126").Value = 3241.6572
$ws.Range("I126").Value = 3502.2856
$ws.Range("J126").Value = 2199.1428
$ws.Range("K126").Value = 10506.8568
$ws.Range("L126").Value = 6597.428400000001
$ws.Range("M126").Value = -8036.856800000001
$ws.Range("N126").Value = -11537.4284

$ws.Range("H134").Value = 43205.766
$ws.Range("I134").Value = 44958.168
$ws.Range("K134").Value = 134874.504
$ws.Range("M134").Value = -132339.504

$ws.Range("H136").Value = 2674.0588
$ws.Range("I136").Value = 2496.3333
$ws.Range("K136").Value = 7488.999899999999
$ws.Range("M136").Value = -4938.999899999999

$ws.Range("H140").Value = 82380.86
$ws.Range("J140").Value = 82380.86
$ws.Range("L140").Value = 82380.86
$ws.Range("N140").Value = -92740.86

$ws.Range("H141").Value = 57163.57
$ws.Range("J141").Value = 57163.57
$ws.Range("L141").Value = 57163.57
$ws.Range("N141").Value = -67523.57000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1854526.4
$ws.Range("J68").Value = 2086072
$ws.Range("L68").Value = 6258216
$ws.Range("N68").Value = -6259838

$ws.Range("H71").Value = 1854526.4
$ws.Range("J71").Value = 2086072
$ws.Range("L71").Value = 18774648
$ws.Range("N71").Value = -18782760

$ws.Range("H122").Value = 1549.6086
$ws.Range("I122").Value = 615.0909
$ws.Range("J122").Value = 2406.25
$ws.Range("K122").Value = 5535.8181
$ws.Range("L122").Value = 21656.25
$ws.Range("M122").Value = -3085.8181
$ws.Range("N122").Value = -26556.25

$ws.Range("H132").Value = 1387.6086
$ws.Range("I132").Value = 1201.2667
$ws.Range("K132").Value = 10811.4003
$ws.Range("M132").Value = -8281.400299999999

$ws.Range("H139").Value = 2752.8076
$ws.Range("J139").Value = 4965
$ws.Range("L139").Value = 14895
$ws.Range("N139").Value = -25175

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1125.8462
$ws.Range("I102").Value = 983.6
$ws.Range("K102").Value = 983.6
$ws.Range("M102").Value = 638.4

$ws.Range("H126").Value = 18843.666
$ws.Range("I126").Value = 27221.77
$ws.Range("J126").Value = 5229.25
$ws.Range("K126").Value = 81665.31
$ws.Range("L126").Value = 15687.75
$ws.Range("M126").Value = -79195.31
$ws.Range("N126").Value = -20627.75

$ws.Range("H132").Value = 19627.55
$ws.Range("I132").Value = 20757.371
$ws.Range("K132").Value = 62272.113
$ws.Range("M132").Value = -59742.113

$ws.Range("H133").Value = 68724.5
$ws.Range("J133").Value = 68724.5
$ws.Range("L133").Value = 68724.5
$ws.Range("N133").Value = -78844.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3498.7144
$ws.Range("I16").Value = 3488.65
$ws.Range("J16").Value = 3700
$ws.Range("K16").Value = 3488.65
$ws.Range("L16").Value = 3700
$ws.Range("M16").Value = -3318.65
$ws.Range("N16").Value = -4040

$ws.Range("H100").Value = 93578.91
$ws.Range("I100").Value = 127171.125
$ws.Range("K100").Value = 127171.125
$ws.Range("M100").Value = -126630.125

$ws.Range("H122").Value = 283188.88
$ws.Range("I122").Value = 420311.78
$ws.Range("K122").Value = 1260935.34
$ws.Range("M122").Value = -1258485.34

$ws.Range("H132").Value = 4104.7896
$ws.Range("I132").Value = 3237.8462
$ws.Range("J132").Value = 5983.1665
$ws.Range("K132").Value = 9713.5386
$ws.Range("L132").Value = 17949.4995
$ws.Range("M132").Value = -7183.5386
$ws.Range("N132").Value = -23009.4995

$ws.Range("H134").Value = 50000
$ws.Range("J134").Value = 50000
$ws.Range("L134").Value = 50000
$ws.Range("N134").Value = -60140

$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2439.4324
$ws.Range("I132").Value = 2440.5278
$ws.Range("K132").Value = 7321.5834
$ws.Range("M132").Value = -4791.5834

$ws.Range("H133").Value = 74875
$ws.Range("J133").Value = 74875
$ws.Range("L133").Value = 74875
$ws.Range("N133").Value = -84995

$ws.Range("H136").Value = 2052.6667
$ws.Range("I136").Value = 1788.4166
$ws.Range("J136").Value = 4166.6665
$ws.Range("K136").Value = 5365.2498
$ws.Range("L136").Value = 12499.9995
$ws.Range("M136").Value = -2815.2498
$ws.Range("N136").Value = -17599.9995

$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
